# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity, and Temperature sheets,
# matching the source system's latest export (rows dated 2026-01-28).

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($ws, $rows, $forceTextValue)

    if ($rows.Count -eq 0) { return }

    $firstRow = $rows[0][0]
    $lastRow = $rows[$rows.Count - 1][0]

    # Column A holds plain-text dates like "2026-01-28"; Excel's smart entry
    # would otherwise silently convert them to date serials, so force the
    # column to Text before writing, then restore the default style.
    $dateRange = $ws.Range("A$firstRow" + ":A$lastRow")
    $dateRange.NumberFormat = "@"

    $valueRange = $null
    if ($forceTextValue) {
        # Column E sometimes holds percentage-looking text (e.g. "87.9%")
        # that Excel would otherwise convert to a numeric percentage.
        $valueRange = $ws.Range("E$firstRow" + ":E$lastRow")
        $valueRange.NumberFormat = "@"
    }

    foreach ($row in $rows) {
        $r = $row[0]
        $ws.Cells.Item($r, 1).Value2 = $row[1]
        $ws.Cells.Item($r, 2).Value2 = $row[2]
        $ws.Cells.Item($r, 3).Value2 = $row[3]
        $ws.Cells.Item($r, 4).Value2 = $row[4]
        $ws.Cells.Item($r, 5).Value2 = $row[5]
        $ws.Cells.Item($r, 6).Value2 = $row[6]
    }

    $dateRange.Style = "Normal"
    if ($forceTextValue) {
        $valueRange.Style = "Normal"
    }
}

# --- PIR sheet: rows 211-223 (No Motion / Inactive) ---
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(211, "2026-01-28", "16:54:39", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(212, "2026-01-28", "16:54:41", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(213, "2026-01-28", "16:54:45", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(214, "2026-01-28", "16:54:49", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(215, "2026-01-28", "16:54:55", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(216, "2026-01-28", "16:54:59", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(217, "2026-01-28", "16:55:05", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(218, "2026-01-28", "16:55:09", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(219, "2026-01-28", "16:55:14", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(220, "2026-01-28", "16:55:19", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(221, "2026-01-28", "16:55:25", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(222, "2026-01-28", "16:55:29", "16:00", "Bathroom", "No Motion", "Inactive"),
    @(223, "2026-01-28", "16:55:35", "16:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows $wsPIR $pirRows $false

# --- Humidity sheet: rows 205-217 (percentage readings) ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(205, "2026-01-28", "16:54:38", "16:00", "Bathroom", "87.9%", "Active"),
    @(206, "2026-01-28", "16:54:39", "16:00", "Bathroom", "87.0%", "Active"),
    @(207, "2026-01-28", "16:54:41", "16:00", "Bathroom", "88.0%", "Active"),
    @(208, "2026-01-28", "16:54:44", "16:00", "Bathroom", "88.0%", "Active"),
    @(209, "2026-01-28", "16:54:48", "16:00", "Bathroom", "87.1%", "Active"),
    @(210, "2026-01-28", "16:54:52", "16:00", "Bathroom", "88.0%", "Active"),
    @(211, "2026-01-28", "16:55:00", "16:00", "Bathroom", "88.0%", "Active"),
    @(212, "2026-01-28", "16:55:04", "16:00", "Bathroom", "88.0%", "Active"),
    @(213, "2026-01-28", "16:55:12", "16:00", "Bathroom", "88.0%", "Active"),
    @(214, "2026-01-28", "16:55:20", "16:00", "Bathroom", "87.1%", "Active"),
    @(215, "2026-01-28", "16:55:24", "16:00", "Bathroom", "88.1%", "Active"),
    @(216, "2026-01-28", "16:55:28", "16:00", "Bathroom", "87.1%", "Active"),
    @(217, "2026-01-28", "16:55:32", "16:00", "Bathroom", "88.0%", "Active")
)
Add-LogRows $wsHumidity $humidityRows $true

# --- Temperature sheet: rows 205-217 (Celsius readings) ---
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(205, "2026-01-28", "16:54:38", "16:00", "Bathroom", "22.8C", "Active"),
    @(206, "2026-01-28", "16:54:40", "16:00", "Bathroom", "22.8C", "Active"),
    @(207, "2026-01-28", "16:54:42", "16:00", "Bathroom", "22.8C", "Active"),
    @(208, "2026-01-28", "16:54:44", "16:00", "Bathroom", "22.8C", "Active"),
    @(209, "2026-01-28", "16:54:48", "16:00", "Bathroom", "22.8C", "Active"),
    @(210, "2026-01-28", "16:54:52", "16:00", "Bathroom", "22.8C", "Active"),
    @(211, "2026-01-28", "16:55:01", "16:00", "Bathroom", "22.8C", "Active"),
    @(212, "2026-01-28", "16:55:04", "16:00", "Bathroom", "22.8C", "Active"),
    @(213, "2026-01-28", "16:55:12", "16:00", "Bathroom", "22.8C", "Active"),
    @(214, "2026-01-28", "16:55:21", "16:00", "Bathroom", "22.8C", "Active"),
    @(215, "2026-01-28", "16:55:24", "16:00", "Bathroom", "22.8C", "Active"),
    @(216, "2026-01-28", "16:55:28", "16:00", "Bathroom", "22.8C", "Active"),
    @(217, "2026-01-28", "16:55:32", "16:00", "Bathroom", "22.8C", "Active")
)
Add-LogRows $wsTemperature $temperatureRows $false
